# Suppress the "permanently delete this sheet" confirmation so the
# automation doesn't stall on the alert dialog.
$excel.DisplayAlerts = $false

$wb = $excel.ActiveWorkbook

# Remove the "Desarquivamentos Pendentes" tab entirely (its data, the
# shared strings it alone used, and the cell styles it alone used all go
# away with it).
$wb.Worksheets("Desarquivamentos Pendentes").Delete() | Out-Null

# Normalize the casing/wording of the remaining tab names.
$wb.Worksheets("Paineis DARQ").Name = "PAINEIS DARQ"
$wb.Worksheets("Recolhimento x Eliminacao").Name = "RECOLHIMENTO X ELIMINAÇÃO"

$excel.DisplayAlerts = $true
